$wb = $excel.ActiveWorkbook

# Target row-2 values per sheet, columns A..O (eb, gb, hp, st, wi, ieh, chp, ac, ab_ct, ab_hp, cp_ct, cp_hp, ttes, btes, ites)
$values = @{
    1 = @(189.0080070514287, 0, 34571.91593949074, 0, 695202.8899878451, 866.090703915558, 0, 2534.277928792126, 0, 0, 0, 0, 0, 218547.9617955566, 1995.76246267978)
    2 = @(189.0080070514287, 0, 148350.0577205477, 0, 695202.8899878451, 866.090703915558, 0, 7410.946080941624, 0, 0, 0, 0, 0, 233137.6946312535, 5738.328092052148)
    3 = @(189.0080070514287, 0, 255308.0505297226, 0, 695202.8899878451, 14355.12988362158, 0, 11593.76417664664, 0, 0, 0, 0, 0, 234698.4240172835, 9655.305390073228)
    4 = @(189.0080070514287, 0, 255308.0505297226, 0, 695202.8899878451, 14355.12988362158, 0, 11593.76417664664, 0, 0, 0, 0, 0, 234925.0693451238, 9655.305390073228)
    5 = @(189.0080070514287, 0, 255308.0505297226, 0, 695202.8899878451, 14355.12988362158, 0, 11593.76417664664, 0, 0, 0, 0, 0, 234925.0693451238, 9655.305390073228)
    6 = @(189.0080070514287, 0, 255308.0505297226, 0, 695202.8899878451, 14355.12988362158, 0, 11593.76417664664, 0, 0, 0, 0, 0, 234925.0693451238, 9655.305390073228)
}

$idx = 0
foreach ($ws in $wb.Worksheets) {
    $idx = $idx + 1

    # Insert a new column before column B (so "gb" is placed right after "eb")
    $ws.Columns.Item(2).Insert()
    # Insert a new column before column N (so "btes" is placed right after "ttes" and before "ites")
    $ws.Columns.Item(14).Insert()

    $ws.Cells.Item(1, 2).Value = "gb"
    $ws.Cells.Item(1, 14).Value = "btes"

    $row = $values[$idx]
    for ($c = 1; $c -le 15; $c++) {
        $ws.Cells.Item(2, $c).Value = $row[$c - 1]
    }
}
